# Update training history metrics (128 dense layers, 20 lstm, 50 epochs, 0.5 dropout run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.05790764093399048
$ws.Cells.Item(2, 2).Value = 0.9890550374984741
$ws.Cells.Item(2, 3).Value = 0.008031287230551243
$ws.Cells.Item(2, 4).Value = 0.9989309906959534

$ws.Cells.Item(3, 1).Value = 0.01109076477587223
$ws.Cells.Item(3, 2).Value = 0.9979838132858276
$ws.Cells.Item(3, 3).Value = 0.005349431652575731
$ws.Cells.Item(3, 4).Value = 0.9989309906959534

$ws.Cells.Item(4, 1).Value = 0.006748031824827194
$ws.Cells.Item(4, 2).Value = 0.9984027743339539
$ws.Cells.Item(4, 3).Value = 0.002753283828496933
$ws.Cells.Item(4, 4).Value = 0.9991776943206787

$ws.Cells.Item(5, 1).Value = 0.003380171256139874
$ws.Cells.Item(5, 2).Value = 0.9989264607429504
$ws.Cells.Item(5, 3).Value = 0.0004931573639623821
$ws.Cells.Item(5, 4).Value = 0.9998355507850647

$ws.Cells.Item(6, 1).Value = 0.002814432606101036
$ws.Cells.Item(6, 2).Value = 0.9992406368255615
$ws.Cells.Item(6, 3).Value = 0.0001125216076616198
$ws.Cells.Item(6, 4).Value = 1

$ws.Cells.Item(7, 1).Value = 0.001504226471297443
$ws.Cells.Item(7, 2).Value = 0.999607264995575
$ws.Cells.Item(7, 3).Value = 0.0001134235862991773
$ws.Cells.Item(7, 4).Value = 1

$ws.Cells.Item(8, 1).Value = 0.001523348968476057
$ws.Cells.Item(8, 2).Value = 0.999450147151947
$ws.Cells.Item(8, 3).Value = 0.0001273883681278676
$ws.Cells.Item(8, 4).Value = 1

$ws.Cells.Item(9, 1).Value = 0.001455856254324317
$ws.Cells.Item(9, 2).Value = 0.9995810389518738
$ws.Cells.Item(9, 3).Value = 0.00007454854494426399
$ws.Cells.Item(9, 4).Value = 1

$ws.Cells.Item(10, 1).Value = 0.0007962927338667214
$ws.Cells.Item(10, 2).Value = 0.9998428821563721
$ws.Cells.Item(10, 3).Value = 0.0001758038124535233
$ws.Cells.Item(10, 4).Value = 0.99991774559021

$ws.Cells.Item(11, 1).Value = 0.0007505666580982506
$ws.Cells.Item(11, 2).Value = 0.9998167157173157
$ws.Cells.Item(11, 3).Value = 0.0000503523915540427
$ws.Cells.Item(11, 4).Value = 1

$ws.Cells.Item(12, 1).Value = 0.0008993609808385372
$ws.Cells.Item(12, 2).Value = 0.9997119903564453
$ws.Cells.Item(12, 3).Value = 0.0001186896115541458
$ws.Cells.Item(12, 4).Value = 1

$ws.Cells.Item(13, 1).Value = 0.00153013167437166
$ws.Cells.Item(13, 2).Value = 0.9996334314346313
$ws.Cells.Item(13, 3).Value = 0.0002088810579152778
$ws.Cells.Item(13, 4).Value = 0.99991774559021

$ws.Cells.Item(14, 1).Value = 0.0007439260953105986
$ws.Cells.Item(14, 2).Value = 0.9997643232345581
$ws.Cells.Item(14, 3).Value = 0.00005022064578952268
$ws.Cells.Item(14, 4).Value = 1

$ws.Cells.Item(15, 1).Value = 0.0006056775455363095
$ws.Cells.Item(15, 2).Value = 0.9997905492782593
$ws.Cells.Item(15, 3).Value = 0.00003790808841586113
$ws.Cells.Item(15, 4).Value = 1

$ws.Cells.Item(16, 1).Value = 0.0003594601294025779
$ws.Cells.Item(16, 2).Value = 0.9998691082000732
$ws.Cells.Item(16, 3).Value = 0.0002261483896290883
$ws.Cells.Item(16, 4).Value = 0.99991774559021

$ws.Cells.Item(17, 1).Value = 0.0009670010185800493
$ws.Cells.Item(17, 2).Value = 0.9997643232345581
$ws.Cells.Item(17, 3).Value = 0.00001421861634298693
$ws.Cells.Item(17, 4).Value = 1

$ws.Cells.Item(18, 1).Value = 0.0003195509780198336
$ws.Cells.Item(18, 2).Value = 0.9998428821563721
$ws.Cells.Item(18, 3).Value = 0.00003950777681893669
$ws.Cells.Item(18, 4).Value = 1

$ws.Cells.Item(19, 1).Value = 0.0003447523340582848
$ws.Cells.Item(19, 2).Value = 0.9998952746391296
$ws.Cells.Item(19, 3).Value = 0.000008705368600203656
$ws.Cells.Item(19, 4).Value = 1

$ws.Cells.Item(20, 1).Value = 0.0008060195832513273
$ws.Cells.Item(20, 2).Value = 0.9998167157173157
$ws.Cells.Item(20, 3).Value = 0.0002927497553173453
$ws.Cells.Item(20, 4).Value = 0.99991774559021

$ws.Cells.Item(21, 1).Value = 0.0001099277869798243
$ws.Cells.Item(21, 2).Value = 0.9999476075172424
$ws.Cells.Item(21, 3).Value = 0.00005061726187705062
$ws.Cells.Item(21, 4).Value = 1

$ws.Cells.Item(22, 1).Value = 0.0005496108205989003
$ws.Cells.Item(22, 2).Value = 0.9998428821563721
$ws.Cells.Item(22, 3).Value = 0.0001494950556661934
$ws.Cells.Item(22, 4).Value = 0.99991774559021

$ws.Cells.Item(23, 1).Value = 0.0003722425026353449
$ws.Cells.Item(23, 2).Value = 0.9998691082000732
$ws.Cells.Item(23, 3).Value = 0.00001518455519544659
$ws.Cells.Item(23, 4).Value = 1

$ws.Cells.Item(24, 1).Value = 0.0003220294311176986
$ws.Cells.Item(24, 2).Value = 0.999921441078186
$ws.Cells.Item(24, 3).Value = 0.000002944539573945804
$ws.Cells.Item(24, 4).Value = 1

$ws.Cells.Item(25, 1).Value = 0.0001521692465757951
$ws.Cells.Item(25, 2).Value = 0.9999476075172424
$ws.Cells.Item(25, 3).Value = 0.00001343099847872509
$ws.Cells.Item(25, 4).Value = 1

$ws.Cells.Item(26, 1).Value = 0.0002951617643702775
$ws.Cells.Item(26, 2).Value = 0.9999476075172424
$ws.Cells.Item(26, 3).Value = 0.000007095762612152612
$ws.Cells.Item(26, 4).Value = 1

$ws.Cells.Item(27, 1).Value = 0.0001217396347783506
$ws.Cells.Item(27, 2).Value = 0.9999476075172424
$ws.Cells.Item(27, 3).Value = 0.0000001185730695851817
$ws.Cells.Item(27, 4).Value = 1

$ws.Cells.Item(28, 1).Value = 0.0003602619399316609
$ws.Cells.Item(28, 2).Value = 0.9999476075172424
$ws.Cells.Item(28, 3).Value = 0.00003858376294374466
$ws.Cells.Item(28, 4).Value = 1

$ws.Cells.Item(29, 1).Value = 0.0007607684237882495
$ws.Cells.Item(29, 2).Value = 0.9997381567955017
$ws.Cells.Item(29, 3).Value = 0.000001602393240318634
$ws.Cells.Item(29, 4).Value = 1

$ws.Cells.Item(30, 1).Value = 0.0007169448072090745
$ws.Cells.Item(30, 2).Value = 0.9998167157173157
$ws.Cells.Item(30, 3).Value = 0.000004956886641593883
$ws.Cells.Item(30, 4).Value = 1

$ws.Cells.Item(31, 1).Value = 0.0002058119425782934
$ws.Cells.Item(31, 2).Value = 0.999921441078186
$ws.Cells.Item(31, 3).Value = 0.0000432564556831494
$ws.Cells.Item(31, 4).Value = 1

$ws.Cells.Item(32, 1).Value = 0.0002373325551161543
$ws.Cells.Item(32, 2).Value = 0.9999738335609436
$ws.Cells.Item(32, 3).Value = 0.000004065311713929987
$ws.Cells.Item(32, 4).Value = 1

$ws.Cells.Item(33, 1).Value = 0.00007145714334910735
$ws.Cells.Item(33, 2).Value = 0.9999738335609436
$ws.Cells.Item(33, 3).Value = 0.000103184866020456
$ws.Cells.Item(33, 4).Value = 0.99991774559021

$ws.Cells.Item(34, 1).Value = 0.0005832509486936033
$ws.Cells.Item(34, 2).Value = 0.999921441078186
$ws.Cells.Item(34, 3).Value = 0.0000008074614470388042
$ws.Cells.Item(34, 4).Value = 1

$ws.Cells.Item(35, 1).Value = 0.0005189369549043477
$ws.Cells.Item(35, 2).Value = 0.9998952746391296
$ws.Cells.Item(35, 3).Value = 0.00004281098154024221
$ws.Cells.Item(35, 4).Value = 1

$ws.Cells.Item(36, 1).Value = 0.0003146185190416873
$ws.Cells.Item(36, 2).Value = 0.9998952746391296
$ws.Cells.Item(36, 3).Value = 0.006118983961641788
$ws.Cells.Item(36, 4).Value = 0.9978620409965515

$ws.Cells.Item(37, 1).Value = 0.0003891826490871608
$ws.Cells.Item(37, 2).Value = 0.999921441078186
$ws.Cells.Item(37, 3).Value = 0.00000005103469646883241
$ws.Cells.Item(37, 4).Value = 1

$ws.Cells.Item(38, 1).Value = 0.00001863564102677628
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(38, 3).Value = 0.00000001679110361862968
$ws.Cells.Item(38, 4).Value = 1

$ws.Cells.Item(39, 1).Value = 0.0003183637745678425
$ws.Cells.Item(39, 2).Value = 0.9999476075172424
$ws.Cells.Item(39, 3).Value = 0.00000006266787977438071
$ws.Cells.Item(39, 4).Value = 1

$ws.Cells.Item(40, 1).Value = 0.00008432302274741232
$ws.Cells.Item(40, 2).Value = 0.9999476075172424
$ws.Cells.Item(40, 3).Value = 0.000000002548657951351174
$ws.Cells.Item(40, 4).Value = 1

$ws.Cells.Item(41, 1).Value = 0.0004737147246487439
$ws.Cells.Item(41, 2).Value = 0.999921441078186
$ws.Cells.Item(41, 3).Value = 0.000000001019467399387963
$ws.Cells.Item(41, 4).Value = 1

$ws.Cells.Item(42, 1).Value = 0.000375552917830646
$ws.Cells.Item(42, 2).Value = 0.9998952746391296
$ws.Cells.Item(42, 3).Value = 0.000000001911499625251167
$ws.Cells.Item(42, 4).Value = 1

$ws.Cells.Item(43, 1).Value = 0.0009642954682931304
$ws.Cells.Item(43, 2).Value = 0.9997381567955017
$ws.Cells.Item(43, 3).Value = 0.00000000290155899307365
$ws.Cells.Item(43, 4).Value = 1

$ws.Cells.Item(44, 1).Value = 0.0001914624590426683
$ws.Cells.Item(44, 2).Value = 0.9999476075172424
$ws.Cells.Item(44, 3).Value = 0.000000001509595892734694
$ws.Cells.Item(44, 4).Value = 1

$ws.Cells.Item(45, 1).Value = 0.0003856393159367144
$ws.Cells.Item(45, 2).Value = 0.999921441078186
$ws.Cells.Item(45, 3).Value = 0.000001704424448689679
$ws.Cells.Item(45, 4).Value = 1

$ws.Cells.Item(46, 1).Value = 0.0005355846369639039
$ws.Cells.Item(46, 2).Value = 0.9998952746391296
$ws.Cells.Item(46, 3).Value = 0.00000002116303399191111
$ws.Cells.Item(46, 4).Value = 1

$ws.Cells.Item(47, 1).Value = 0.0002260941255372018
$ws.Cells.Item(47, 2).Value = 0.9999738335609436
$ws.Cells.Item(47, 3).Value = 0.00000009343176543552545
$ws.Cells.Item(47, 4).Value = 1

$ws.Cells.Item(48, 1).Value = 0.000163688775501214
$ws.Cells.Item(48, 2).Value = 0.9999738335609436
$ws.Cells.Item(48, 3).Value = 0.0000001197692398591244
$ws.Cells.Item(48, 4).Value = 1

$ws.Cells.Item(49, 1).Value = 0.0002383910905336961
$ws.Cells.Item(49, 2).Value = 0.9999476075172424
$ws.Cells.Item(49, 3).Value = 0.00000001411557892794235
$ws.Cells.Item(49, 4).Value = 1

$ws.Cells.Item(50, 1).Value = 0.0001218413381138816
$ws.Cells.Item(50, 2).Value = 0.9999738335609436
$ws.Cells.Item(50, 3).Value = 0.000000006244146533163075
$ws.Cells.Item(50, 4).Value = 1

$ws.Cells.Item(51, 1).Value = 0.00005417723514256068
$ws.Cells.Item(51, 2).Value = 0.9999738335609436
$ws.Cells.Item(51, 3).Value = 0.000000001529199211702803
$ws.Cells.Item(51, 4).Value = 1

